$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 14.956293
$ws.Range("H2").Value = 44.868879
$ws.Range("I2").Value = 0.195346866793292
$ws.Range("J2").Value = 0.1953468667932919
$ws.Range("M2").Value = 73.202511
$ws.Range("N2").Value = 219.607533
$ws.Range("O2").Value = 0.3264904632507938
$ws.Range("P2").Value = 0.3264904632507938
$ws.Range("Q2").Value = 1094.838202851723
$ws.Range("R2").Value = 9853.543825665507
$ws.Range("S2").Value = 0.063778889033933
$ws.Range("T2").Value = 0.063778889033933
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 14.956293
$ws.Range("H3").Value = 44.868879
$ws.Range("I3").Value = 0.195346866793292
$ws.Range("J3").Value = 0.1953468667932919
$ws.Range("O3").Value = 0.4449719839907295
$ws.Range("P3").Value = 0.4449719839907295
$ws.Range("Q3").Value = 1492.14872134735
$ws.Range("R3").Value = 13429.33849212615
$ws.Range("S3").Value = 0.08692388288338389
$ws.Range("T3").Value = 0.08692388288338387
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 14.956293
$ws.Range("H4").Value = 44.868879
$ws.Range("I4").Value = 0.195346866793292
$ws.Range("J4").Value = 0.1953468667932919
$ws.Range("M4").Value = 39.54025133333334
$ws.Range("N4").Value = 118.620754
$ws.Range("O4").Value = 0.1763534446908907
$ws.Range("P4").Value = 0.1763534446908907
$ws.Range("Q4").Value = 591.3755842349741
$ws.Range("R4").Value = 5322.380258114767
$ws.Range("S4").Value = 0.03445009286856961
$ws.Range("T4").Value = 0.03445009286856961
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.956293
$ws.Range("H5").Value = 44.868879
$ws.Range("I5").Value = 0.195346866793292
$ws.Range("J5").Value = 0.1953468667932919
$ws.Range("M5").Value = 11.70021233333333
$ws.Range("N5").Value = 35.100637
$ws.Range("O5").Value = 0.05218410806758597
$ws.Range("P5").Value = 0.05218410806758598
$ws.Range("Q5").Value = 174.991803819547
$ws.Range("R5").Value = 1574.926234375923
$ws.Range("S5").Value = 0.01019400200740547
$ws.Range("T5").Value = 0.01019400200740547
$ws.Range("I6").Value = 0.653630228225219
$ws.Range("J6").Value = 0.653630228225219
$ws.Range("M6").Value = 73.202511
$ws.Range("N6").Value = 219.607533
$ws.Range("O6").Value = 0.3264904632507938
$ws.Range("P6").Value = 0.3264904632507938
$ws.Range("Q6").Value = 3663.32645179766
$ws.Range("R6").Value = 32969.93806617894
$ws.Range("S6").Value = 0.2134040360079738
$ws.Range("T6").Value = 0.2134040360079739
$ws.Range("I7").Value = 0.653630228225219
$ws.Range("J7").Value = 0.653630228225219
$ws.Range("O7").Value = 0.4449719839907295
$ws.Range("P7").Value = 0.4449719839907295
$ws.Range("S7").Value = 0.290847139449689
$ws.Range("T7").Value = 0.290847139449689
$ws.Range("I8").Value = 0.653630228225219
$ws.Range("J8").Value = 0.653630228225219
$ws.Range("M8").Value = 39.54025133333334
$ws.Range("N8").Value = 118.620754
$ws.Range("O8").Value = 0.1763534446908907
$ws.Range("P8").Value = 0.1763534446908907
$ws.Range("Q8").Value = 1978.741530056636
$ws.Range("R8").Value = 17808.67377050972
$ws.Range("S8").Value = 0.1152699423016104
$ws.Range("T8").Value = 0.1152699423016104
$ws.Range("I9").Value = 0.653630228225219
$ws.Range("J9").Value = 0.653630228225219
$ws.Range("M9").Value = 11.70021233333333
$ws.Range("N9").Value = 35.100637
$ws.Range("O9").Value = 0.05218410806758597
$ws.Range("P9").Value = 0.05218410806758598
$ws.Range("Q9").Value = 585.5222279512955
$ws.Range("R9").Value = 5269.70005156166
$ws.Range("S9").Value = 0.03410911046594571
$ws.Range("T9").Value = 0.03410911046594572
$ws.Range("G10").Value = 10.01531766666667
$ws.Range("H10").Value = 30.045953
$ws.Range("I10").Value = 0.1308118880876991
$ws.Range("J10").Value = 0.1308118880876991
$ws.Range("M10").Value = 73.202511
$ws.Range("N10").Value = 219.607533
$ws.Range("O10").Value = 0.3264904632507938
$ws.Range("P10").Value = 0.3264904632507938
$ws.Range("Q10").Value = 733.146401662661
$ws.Range("R10").Value = 6598.317614963948
$ws.Range("S10").Value = 0.04270883394046386
$ws.Range("T10").Value = 0.04270883394046387
$ws.Range("G11").Value = 10.01531766666667
$ws.Range("H11").Value = 30.045953
$ws.Range("I11").Value = 0.1308118880876991
$ws.Range("J11").Value = 0.1308118880876991
$ws.Range("O11").Value = 0.4449719839907295
$ws.Range("P11").Value = 0.4449719839907295
$ws.Range("Q11").Value = 999.2010353236723
$ws.Range("R11").Value = 8992.80931791305
$ws.Range("S11").Value = 0.05820762537195673
$ws.Range("T11").Value = 0.05820762537195673
$ws.Range("G12").Value = 10.01531766666667
$ws.Range("H12").Value = 30.045953
$ws.Range("I12").Value = 0.1308118880876991
$ws.Range("J12").Value = 0.1308118880876991
$ws.Range("M12").Value = 39.54025133333334
$ws.Range("N12").Value = 118.620754
$ws.Range("O12").Value = 0.1763534446908907
$ws.Range("P12").Value = 0.1763534446908907
$ws.Range("Q12").Value = 396.0081777231736
$ws.Range("R12").Value = 3564.073599508562
$ws.Range("S12").Value = 0.02306912707078502
$ws.Range("T12").Value = 0.02306912707078502
$ws.Range("G13").Value = 10.01531766666667
$ws.Range("H13").Value = 30.045953
$ws.Range("I13").Value = 0.1308118880876991
$ws.Range("J13").Value = 0.1308118880876991
$ws.Range("M13").Value = 11.70021233333333
$ws.Range("N13").Value = 35.100637
$ws.Range("O13").Value = 0.05218410806758597
$ws.Range("P13").Value = 0.05218410806758598
$ws.Range("Q13").Value = 117.1813432857845
$ws.Range("R13").Value = 1054.632089572061
$ws.Range("S13").Value = 0.00682630170449345
$ws.Range("T13").Value = 0.00682630170449345
$ws.Range("G14").Value = 1.547411
$ws.Range("H14").Value = 4.642232999999999
$ws.Range("I14").Value = 0.0202110168937901
$ws.Range("J14").Value = 0.0202110168937901
$ws.Range("M14").Value = 73.202511
$ws.Range("N14").Value = 219.607533
$ws.Range("O14").Value = 0.3264904632507938
$ws.Range("P14").Value = 0.3264904632507938
$ws.Range("Q14").Value = 113.274370749021
$ws.Range("R14").Value = 1019.469336741189
$ws.Range("S14").Value = 0.00659870426842315
$ws.Range("T14").Value = 0.00659870426842315
$ws.Range("G15").Value = 1.547411
$ws.Range("H15").Value = 4.642232999999999
$ws.Range("I15").Value = 0.0202110168937901
$ws.Range("J15").Value = 0.0202110168937901
$ws.Range("O15").Value = 0.4449719839907295
$ws.Range("P15").Value = 0.4449719839907295
$ws.Range("Q15").Value = 154.3809916701166
$ws.Range("R15").Value = 1389.42892503105
$ws.Range("S15").Value = 0.008993336285699934
$ws.Range("T15").Value = 0.008993336285699932
$ws.Range("G16").Value = 1.547411
$ws.Range("H16").Value = 4.642232999999999
$ws.Range("I16").Value = 0.0202110168937901
$ws.Range("J16").Value = 0.0202110168937901
$ws.Range("M16").Value = 39.54025133333334
$ws.Range("N16").Value = 118.620754
$ws.Range("O16").Value = 0.1763534446908907
$ws.Range("P16").Value = 0.1763534446908907
$ws.Range("Q16").Value = 61.18501985596466
$ws.Range("R16").Value = 550.6651787036819
$ws.Range("S16").Value = 0.003564282449925671
$ws.Range("T16").Value = 0.00356428244992567
$ws.Range("G17").Value = 1.547411
$ws.Range("H17").Value = 4.642232999999999
$ws.Range("I17").Value = 0.0202110168937901
$ws.Range("J17").Value = 0.0202110168937901
$ws.Range("M17").Value = 11.70021233333333
$ws.Range("N17").Value = 35.100637
$ws.Range("O17").Value = 0.05218410806758597
$ws.Range("P17").Value = 0.05218410806758598
$ws.Range("Q17").Value = 18.10503726693566
$ws.Range("R17").Value = 162.945335402421
$ws.Range("S17").Value = 0.001054693889741349
$ws.Range("T17").Value = 0.001054693889741349
